$wb = $excel.ActiveWorkbook

# This script applies updated market-board derived profit metrics
# (currentAveragePrice/NQ/HQ, LevePriceNQ/HQ, LeveProfitNQ/HQ columns H-N)
# to specific rows across the ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets,
# matching a scheduled refresh of the source pricing data.

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H6").Value = 653.1667
$ws.Range("I6").Value = 653.1667
$ws.Range("K6").Value = 1959.5001
$ws.Range("M6").Value = -1847.5001
$ws.Range("H17").Value = 2387.6
$ws.Range("J17").Value = 2387.6
$ws.Range("L17").Value = 7162.799999999999
$ws.Range("N17").Value = -7498.799999999999
$ws.Range("H43").Value = 2420.2222
$ws.Range("I43").Value = 2029.3334
$ws.Range("J43").Value = 2615.6667
$ws.Range("K43").Value = 2029.3334
$ws.Range("L43").Value = 2615.6667
$ws.Range("M43").Value = -1960.3334
$ws.Range("N43").Value = -2753.6667
$ws.Range("H51").Value = 7755.4546
$ws.Range("I51").Value = 8968.333000000001
$ws.Range("K51").Value = 8968.333000000001
$ws.Range("M51").Value = -8484.333000000001
$ws.Range("H74").Value = 5460.5557
$ws.Range("I74").Value = 5393.75
$ws.Range("K74").Value = 5393.75
$ws.Range("M74").Value = -4457.75
$ws.Range("H77").Value = 5460.5557
$ws.Range("I77").Value = 5393.75
$ws.Range("K77").Value = 26968.75
$ws.Range("M77").Value = -22288.75
$ws.Range("H88").Value = 11620
$ws.Range("I88").Value = 13049.5
$ws.Range("K88").Value = 13049.5
$ws.Range("M88").Value = -12643.5
$ws.Range("H91").Value = 11620
$ws.Range("I91").Value = 13049.5
$ws.Range("K91").Value = 13049.5
$ws.Range("M91").Value = -11645.5
$ws.Range("H96").Value = 1533
$ws.Range("I96").Value = 199.5
$ws.Range("J96").Value = 2199.75
$ws.Range("K96").Value = 598.5
$ws.Range("L96").Value = 6599.25
$ws.Range("M96").Value = 774.5
$ws.Range("N96").Value = -9345.25
$ws.Range("H106").Value = 8292
$ws.Range("I106").Value = 6222.1
$ws.Range("J106").Value = 15191.667
$ws.Range("K106").Value = 6222.1
$ws.Range("L106").Value = 15191.667
$ws.Range("M106").Value = -5591.1
$ws.Range("N106").Value = -16453.667
$ws.Range("H107").Value = 2306.1304
$ws.Range("J107").Value = 2399.2856
$ws.Range("L107").Value = 2399.2856
$ws.Range("N107").Value = -6239.2856
$ws.Range("H125").Value = 1560
$ws.Range("I125").Value = 998
$ws.Range("K125").Value = 8982
$ws.Range("M125").Value = -6522

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 2385.6428
$ws.Range("I2").Value = 1949.9166
$ws.Range("K2").Value = 1949.9166
$ws.Range("M2").Value = -1836.9166
$ws.Range("H32").Value = 22675.537
$ws.Range("I32").Value = 22528.441
$ws.Range("K32").Value = 22528.441
$ws.Range("M32").Value = -22241.441
$ws.Range("H50").Value = 5234.5
$ws.Range("I50").Value = 48
$ws.Range("J50").Value = 6963.3335
$ws.Range("K50").Value = 48
$ws.Range("L50").Value = 6963.3335
$ws.Range("M50").Value = 666
$ws.Range("N50").Value = -8391.333500000001
$ws.Range("H116").Value = 2385.6428
$ws.Range("I116").Value = 1949.9166
$ws.Range("K116").Value = 1949.9166
$ws.Range("M116").Value = 344.0834
$ws.Range("H132").Value = 37017.07
$ws.Range("I132").Value = 42371.8
$ws.Range("K132").Value = 127115.4
$ws.Range("M132").Value = -124585.4

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 2385.6428
$ws.Range("I3").Value = 1949.9166
$ws.Range("K3").Value = 1949.9166
$ws.Range("M3").Value = -1835.9166
$ws.Range("H86").Value = 2472.6924
$ws.Range("I86").Value = 2371.6
$ws.Range("J86").Value = 5000
$ws.Range("K86").Value = 2371.6
$ws.Range("L86").Value = 5000
$ws.Range("M86").Value = -1248.6
$ws.Range("N86").Value = -7246
$ws.Range("H89").Value = 2472.6924
$ws.Range("I89").Value = 2371.6
$ws.Range("J89").Value = 5000
$ws.Range("K89").Value = 11858
$ws.Range("L89").Value = 25000
$ws.Range("M89").Value = -6242
$ws.Range("N89").Value = -36232
$ws.Range("H99").Value = 36599.934
$ws.Range("I99").Value = 74627.71000000001
$ws.Range("J99").Value = 3325.625
$ws.Range("K99").Value = 74627.71000000001
$ws.Range("L99").Value = 3325.625
$ws.Range("M99").Value = -73129.71000000001
$ws.Range("N99").Value = -6321.625
$ws.Range("H134").Value = 1584.4728
$ws.Range("I134").Value = 1262.6735
$ws.Range("K134").Value = 3788.020500000001
$ws.Range("M134").Value = -1253.020500000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 817
$ws.Range("I16").Value = 772.13336
$ws.Range("K16").Value = 772.13336
$ws.Range("M16").Value = -485.13336
$ws.Range("H31").Value = 2056.0645
$ws.Range("I31").Value = 2000.3334
$ws.Range("J31").Value = 2432.25
$ws.Range("K31").Value = 2000.3334
$ws.Range("L31").Value = 2432.25
$ws.Range("M31").Value = -1705.3334
$ws.Range("N31").Value = -3022.25
$ws.Range("H34").Value = 2056.0645
$ws.Range("I34").Value = 2000.3334
$ws.Range("J34").Value = 2432.25
$ws.Range("K34").Value = 2000.3334
$ws.Range("L34").Value = 2432.25
$ws.Range("M34").Value = -1798.3334
$ws.Range("N34").Value = -2836.25
$ws.Range("H103").Value = 53511.5
$ws.Range("I103").Value = 48015.668
$ws.Range("K103").Value = 48015.668
$ws.Range("M103").Value = -46843.668
$ws.Range("H113").Value = 817
$ws.Range("I113").Value = 772.13336
$ws.Range("K113").Value = 772.13336
$ws.Range("M113").Value = 1397.86664
$ws.Range("H132").Value = 2388.7334
$ws.Range("I132").Value = 2141.3845
$ws.Range("K132").Value = 6424.1535
$ws.Range("M132").Value = -3894.1535
$ws.Range("H134").Value = 44168.75
$ws.Range("I134").Value = 54423.79
$ws.Range("K134").Value = 163271.37
$ws.Range("M134").Value = -160736.37

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H68").Value = 50399.5
$ws.Range("J68").Value = 50399.5
$ws.Range("L68").Value = 151198.5
$ws.Range("N68").Value = -152820.5
$ws.Range("H71").Value = 50399.5
$ws.Range("J71").Value = 50399.5
$ws.Range("L71").Value = 453595.5
$ws.Range("N71").Value = -461707.5
$ws.Range("H98").Value = 2449.875
$ws.Range("I98").Value = 398.66666
$ws.Range("K98").Value = 1195.99998
$ws.Range("M98").Value = 302.0000199999999
$ws.Range("H113").Value = 987.8889
$ws.Range("I113").Value = 1147.5
$ws.Range("K113").Value = 3442.5
$ws.Range("M113").Value = -1272.5
$ws.Range("H128").Value = 121745
$ws.Range("I128").Value = 121745
$ws.Range("K128").Value = 365235
$ws.Range("M128").Value = -360255

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H80").Value = 2456.0278
$ws.Range("I80").Value = 1719.25
$ws.Range("K80").Value = 1719.25
$ws.Range("M80").Value = -721.25
$ws.Range("H83").Value = 2456.0278
$ws.Range("I83").Value = 1719.25
$ws.Range("K83").Value = 8596.25
$ws.Range("M83").Value = -3604.25
$ws.Range("H113").Value = 177991
$ws.Range("I113").Value = 125210.22
$ws.Range("K113").Value = 125210.22
$ws.Range("M113").Value = -123040.22
$ws.Range("H126").Value = 7527.4614
$ws.Range("I126").Value = 7939
$ws.Range("K126").Value = 23817
$ws.Range("M126").Value = -21347
$ws.Range("H132").Value = 69636.734
$ws.Range("I132").Value = 74403.78999999999
$ws.Range("J132").Value = 2898
$ws.Range("K132").Value = 223211.37
$ws.Range("L132").Value = 8694
$ws.Range("M132").Value = -220681.37
$ws.Range("N132").Value = -13754

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 5690.0312
$ws.Range("J46").Value = 3664.3809
$ws.Range("L46").Value = 3664.3809
$ws.Range("N46").Value = -4040.3809
$ws.Range("H132").Value = 40044
$ws.Range("I132").Value = 56308.914
$ws.Range("K132").Value = 168926.742
$ws.Range("M132").Value = -166396.742
$ws.Range("H136").Value = 3594.6
$ws.Range("I136").Value = 3206.7144
$ws.Range("J136").Value = 4499.6665
$ws.Range("K136").Value = 9620.143199999999
$ws.Range("L136").Value = 13498.9995
$ws.Range("M136").Value = -7070.143199999999
$ws.Range("N136").Value = -18598.9995

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H100").Value = 442.46155
$ws.Range("I100").Value = 378
$ws.Range("K100").Value = 756
$ws.Range("M100").Value = -215
$ws.Range("H122").Value = 46063.523
$ws.Range("I122").Value = 2522.85
$ws.Range("J122").Value = 336334.66
$ws.Range("K122").Value = 7568.549999999999
$ws.Range("L122").Value = 1009003.98
$ws.Range("M122").Value = -5118.549999999999
$ws.Range("N122").Value = -1013903.98
$ws.Range("H126").Value = 73984.07000000001
$ws.Range("J126").Value = 25600.8
$ws.Range("L126").Value = 76802.39999999999
$ws.Range("N126").Value = -81742.39999999999
$ws.Range("H132").Value = 29295.676
$ws.Range("I132").Value = 30429.658
$ws.Range("J132").Value = 7750
$ws.Range("K132").Value = 91288.974
$ws.Range("L132").Value = 23250
$ws.Range("M132").Value = -88758.974
$ws.Range("N132").Value = -28310
$ws.Range("H136").Value = 4017.7727
$ws.Range("I136").Value = 3980.524
$ws.Range("J136").Value = 4800
$ws.Range("K136").Value = 11941.572
$ws.Range("L136").Value = 14400
$ws.Range("M136").Value = -9391.572
$ws.Range("N136").Value = -19500
